# Update the "取得日時" (acquisition timestamp) column on the ランサーズ sheet
# from "2025-12-02 12:39:15" to "2025-12-02 12:52:40" for rows 2-10.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

for ($r = 2; $r -le 10; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value2 -eq "2025-12-02 12:39:15") {
        $cell.Value = "2025-12-02 12:52:40"
    }
}
